$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing test contents from B1:B2
$ws.Range("B1:B2").ClearContents()

# Widen column B to fit the content that used to be there
$ws.Columns("B").ColumnWidth = 30.42578125

# Apply date-style number formats to A1 and B1 (as if preparing cells
# for date entries), matching built-in numFmtId 16 ("d-mmm") and
# built-in numFmtId 14 (the locale short-date format, "mm-dd-yy")
$ws.Range("A1").NumberFormat = "d-mmm"
$ws.Range("B1").NumberFormat = "mm-dd-yy"

# Update the active selection to D11
$ws.Range("D11").Select()
